$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.32 = 24996.84 pesos`n✅ 24996.84 pesos = 6.3 = 961.71 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 158.3
$ws2.Range("O10").Value = 3957
$ws2.Range("N12").Value = 3969
$ws2.Range("O12").Value = 152.7
